# Update "想去人数" (F column) counts across the workbook sheets to reflect
# the latest scrape output, per commit "Update gh-pages to output generated
# at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 14935
$wsExpo.Range("F3").Value  = 18722
$wsExpo.Range("F5").Value  = 125
$wsExpo.Range("F15").Value = 207
$wsExpo.Range("F22").Value = 7787
$wsExpo.Range("F28").Value = 5986
$wsExpo.Range("F30").Value = 69
$wsExpo.Range("F36").Value = 41

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 8

# Sheet "全部类型" (All types / combined)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 14935
$wsAll.Range("F3").Value  = 18722
$wsAll.Range("F5").Value  = 125
$wsAll.Range("F15").Value = 207
$wsAll.Range("F23").Value = 7787
$wsAll.Range("F29").Value = 8
$wsAll.Range("F31").Value = 5986
$wsAll.Range("F33").Value = 69
$wsAll.Range("F39").Value = 41
